$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B5").Value = "SingleUseId1"
$ws.Range("C5").Value = "Default"
$ws.Range("D5").Value = "Center"
$ws.Range("E5").Value = "LTR"
$ws.Range("F5").Value = "Nowa Gra"

$ws.Range("B4").Value = "SingleUseId2"
$ws.Range("C4").Value = "Default"
$ws.Range("D4").Value = "Center"
$ws.Range("E4").Value = "LTR"
$ws.Range("F4").Value = "Wyjdź"
